$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert Corequisites/Concurrent/Recommended columns, shift Terms Typically Offered to G ---
$ws.Range('D1').Value = 'Corequisites'
$ws.Range('E1').Value = 'Concurrent'
$ws.Range('F1').Value = 'Recommended'
$ws.Range('G1').Value = 'Terms Typically Offered'

# --- Data rows 2-41: populate new Corequisites(D)/Concurrent(E)/Recommended(F) columns,
# move old "Terms Typically Offered" value (old column D) into G, and strip the
# "Recommended: ..." suffix out of Prerequisites (C) into F where present. ---
# Row 2
$ws.Range('D2').Value = 'NA'
$ws.Range('E2').Value = 'NA'
$ws.Range('F2').Value = 'NA'
$ws.Range('G2').Value = 'F, W, SP'

# Row 3
$ws.Range('D3').Value = 'NA'
$ws.Range('E3').Value = 'NA'
$ws.Range('F3').Value = 'NA'
$ws.Range('G3').Value = 'W'

# Row 4
$ws.Range('D4').Value = 'NA'
$ws.Range('E4').Value = 'NA'
$ws.Range('F4').Value = 'NA'
$ws.Range('G4').Value = 'F, W, SP'

# Row 5
$ws.Range('D5').Value = 'NA'
$ws.Range('E5').Value = 'NA'
$ws.Range('F5').Value = 'NA'
$ws.Range('G5').Value = 'TBD'

# Row 6
$ws.Range('C6').Value = 'Completion of GE Area D1.'
$ws.Range('D6').Value = 'NA'
$ws.Range('E6').Value = 'NA'
$ws.Range('F6').Value = 'ES 112.'
$ws.Range('G6').Value = 'TBD '

# Row 7
$ws.Range('D7').Value = 'NA'
$ws.Range('E7').Value = 'NA'
$ws.Range('F7').Value = 'NA'
$ws.Range('G7').Value = 'F, W, SP'

# Row 8
$ws.Range('D8').Value = 'NA'
$ws.Range('E8').Value = 'NA'
$ws.Range('F8').Value = 'NA'
$ws.Range('G8').Value = 'F, W, SP'

# Row 9
$ws.Range('D9').Value = 'NA'
$ws.Range('E9').Value = 'NA'
$ws.Range('F9').Value = 'NA'
$ws.Range('G9').Value = 'F, W, SP'

# Row 10
$ws.Range('D10').Value = 'NA'
$ws.Range('E10').Value = 'NA'
$ws.Range('F10').Value = 'NA'
$ws.Range('G10').Value = 'F, W, SP'

# Row 11
$ws.Range('D11').Value = 'NA'
$ws.Range('E11').Value = 'NA'
$ws.Range('F11').Value = 'NA'
$ws.Range('G11').Value = 'TBD'

# Row 12
$ws.Range('D12').Value = 'NA'
$ws.Range('E12').Value = 'NA'
$ws.Range('F12').Value = 'NA'
$ws.Range('G12').Value = 'TBD'

# Row 13
$ws.Range('D13').Value = 'NA'
$ws.Range('E13').Value = 'NA'
$ws.Range('F13').Value = 'NA'
$ws.Range('G13').Value = 'F, W, SP'

# Row 14
$ws.Range('D14').Value = 'NA'
$ws.Range('E14').Value = 'NA'
$ws.Range('F14').Value = 'NA'
$ws.Range('G14').Value = 'W, SP'

# Row 15
$ws.Range('D15').Value = 'NA'
$ws.Range('E15').Value = 'NA'
$ws.Range('F15').Value = 'NA'
$ws.Range('G15').Value = 'SP'

# Row 16
$ws.Range('D16').Value = 'NA'
$ws.Range('E16').Value = 'NA'
$ws.Range('F16').Value = 'NA'
$ws.Range('G16').Value = 'TBD'

# Row 17
$ws.Range('C17').Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D.'
$ws.Range('D17').Value = 'NA'
$ws.Range('E17').Value = 'NA'
$ws.Range('F17').Value = 'ES 112 (GE Area D1) or any ES course in GE Area D3.'
$ws.Range('G17').Value = 'W '

# Row 18
$ws.Range('D18').Value = 'NA'
$ws.Range('E18').Value = 'NA'
$ws.Range('F18').Value = 'NA'
$ws.Range('G18').Value = 'SP'

# Row 19
$ws.Range('C19').Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D.'
$ws.Range('D19').Value = 'NA'
$ws.Range('E19').Value = 'NA'
$ws.Range('F19').Value = 'ES 112 (GE Area D1) or any ES course in GE Area D3.'
$ws.Range('G19').Value = 'F, W, SP '

# Row 20
$ws.Range('C20').Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D.'
$ws.Range('D20').Value = 'NA'
$ws.Range('E20').Value = 'NA'
$ws.Range('F20').Value = 'ES 112 (GE Area D1) or any ES course in GE Area D3.'
$ws.Range('G20').Value = 'F, W, SP '

# Row 21
$ws.Range('C21').Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D.'
$ws.Range('D21').Value = 'NA'
$ws.Range('E21').Value = 'NA'
$ws.Range('F21').Value = 'ES 112 (GE Area D1) or any ES course in GE Area D3.'
$ws.Range('G21').Value = 'F, W, SP '

# Row 22
$ws.Range('C22').Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D.'
$ws.Range('D22').Value = 'NA'
$ws.Range('E22').Value = 'NA'
$ws.Range('F22').Value = 'ES 112 (GE Area D1) or any ES course in GE Area D3.'
$ws.Range('G22').Value = 'F, W, SP '

# Row 23
$ws.Range('D23').Value = 'NA'
$ws.Range('E23').Value = 'NA'
$ws.Range('F23').Value = 'NA'
$ws.Range('G23').Value = 'SP'

# Row 24
$ws.Range('C24').Value = 'Completion of a course in GE Area D1 or D3.'
$ws.Range('D24').Value = 'NA'
$ws.Range('E24').Value = 'NA'
$ws.Range('F24').Value = 'ES 112 or ES 212.'
$ws.Range('G24').Value = 'TBD '

# Row 25
$ws.Range('D25').Value = 'NA'
$ws.Range('E25').Value = 'NA'
$ws.Range('F25').Value = 'NA'
$ws.Range('G25').Value = 'W'

# Row 26
$ws.Range('C26').Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D.'
$ws.Range('D26').Value = 'NA'
$ws.Range('E26').Value = 'NA'
$ws.Range('F26').Value = 'ES 112 (GE Area D1) or any ES course in GE Area D3.'
$ws.Range('G26').Value = 'TBD '

# Row 27
$ws.Range('C27').Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D.'
$ws.Range('D27').Value = 'NA'
$ws.Range('E27').Value = 'NA'
$ws.Range('F27').Value = 'ES 112 (GE Area D1) or any ES course in GE Area D3.'
$ws.Range('G27').Value = 'W '

# Row 28
$ws.Range('C28').Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area C.'
$ws.Range('D28').Value = 'NA'
$ws.Range('E28').Value = 'NA'
$ws.Range('F28').Value = 'Completion of an Ethnic Studies (ES) course.'
$ws.Range('G28').Value = 'F, W, SP '

# Row 29
$ws.Range('D29').Value = 'NA'
$ws.Range('E29').Value = 'NA'
$ws.Range('F29').Value = 'NA'
$ws.Range('G29').Value = 'F, W, SP'

# Row 30
$ws.Range('D30').Value = 'NA'
$ws.Range('E30').Value = 'NA'
$ws.Range('F30').Value = 'NA'
$ws.Range('G30').Value = 'F, W, SP'

# Row 31
$ws.Range('D31').Value = 'NA'
$ws.Range('E31').Value = 'NA'
$ws.Range('F31').Value = 'NA'
$ws.Range('G31').Value = 'SP'

# Row 32
$ws.Range('C32').Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area C.'
$ws.Range('D32').Value = 'NA'
$ws.Range('E32').Value = 'NA'
$ws.Range('F32').Value = 'Lower-division Ethnic Studies (ES) course and an introductory natural resources course.'
$ws.Range('G32').Value = 'W '

# Row 33
$ws.Range('D33').Value = 'NA'
$ws.Range('E33').Value = 'NA'
$ws.Range('F33').Value = 'NA'
$ws.Range('G33').Value = 'W, SP'

# Row 34
$ws.Range('C34').Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D.'
$ws.Range('D34').Value = 'NA'
$ws.Range('E34').Value = 'NA'
$ws.Range('F34').Value = 'ES 112 (GE Area D1) or any ES course in GE Area D3.'
$ws.Range('G34').Value = 'SP '

# Row 35
$ws.Range('D35').Value = 'NA'
$ws.Range('E35').Value = 'NA'
$ws.Range('F35').Value = 'NA'
$ws.Range('G35').Value = 'F'

# Row 36
$ws.Range('D36').Value = 'NA'
$ws.Range('E36').Value = 'NA'
$ws.Range('F36').Value = 'NA'
$ws.Range('G36').Value = 'F, W, SP'

# Row 37
$ws.Range('D37').Value = 'NA'
$ws.Range('E37').Value = 'NA'
$ws.Range('F37').Value = 'NA'
$ws.Range('G37').Value = 'TBD'

# Row 38
$ws.Range('D38').Value = 'NA'
$ws.Range('E38').Value = 'NA'
$ws.Range('F38').Value = 'NA'
$ws.Range('G38').Value = 'TBD'

# Row 39
$ws.Range('D39').Value = 'NA'
$ws.Range('E39').Value = 'NA'
$ws.Range('F39').Value = 'NA'
$ws.Range('G39').Value = 'W'

# Row 40
$ws.Range('D40').Value = 'NA'
$ws.Range('E40').Value = 'NA'
$ws.Range('F40').Value = 'NA'
$ws.Range('G40').Value = 'SP'

# Row 41
$ws.Range('D41').Value = 'NA'
$ws.Range('E41').Value = 'NA'
$ws.Range('F41').Value = 'NA'
$ws.Range('G41').Value = 'TBD'

